$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column F ("Voltage"), shifting
# Voltage..Library Path one column to the right (F..N -> G..O).
$ws.Columns.Item(6).Insert()

# Give the new column a header and a width similar to the other data columns.
$ws.Cells.Item(1, 6).Value = "Dielectric"
$ws.Columns.Item(6).ColumnWidth = 9.5703125

# Populate the Dielectric value for every data row, based on the dielectric
# type already present in each row's Description (column B).
$dielectrics = @("C0G","C0G","C0G","C0G","X7R","X7R","X7R","C0G","X7R","X7R","X7R","X7R","X7R","X7R","X7R","X7R","X7R","X7R","X5R","X5R","X5R","X5R","X5R","X5R")

for ($i = 0; $i -lt $dielectrics.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $dielectrics[$i]
}
